# Auto-generated Excel COM-interop script
# Applies market-data value refresh to columns H-N across all 8 class sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1292.25
$ws.Range("I43").Value = 977.25
$ws.Range("J43").Value = 1449.75
$ws.Range("K43").Value = 977.25
$ws.Range("L43").Value = 1449.75
$ws.Range("M43").Value = -908.25
$ws.Range("N43").Value = -1587.75
$ws.Range("H62").Value = 2616.4736
$ws.Range("J62").Value = 2583.3333
$ws.Range("L62").Value = 2583.3333
$ws.Range("N62").Value = -3831.3333
$ws.Range("H65").Value = 2616.4736
$ws.Range("J65").Value = 2583.3333
$ws.Range("L65").Value = 12916.6665
$ws.Range("N65").Value = -19156.6665
$ws.Range("H105").Value = 45000
$ws.Range("J105").Value = 45000
$ws.Range("L105").Value = 45000
$ws.Range("N105").Value = -51988
$ws.Range("H112").Value = 1398.7838
$ws.Range("J112").Value = 1522.8788
$ws.Range("L112").Value = 4568.636399999999
$ws.Range("N112").Value = -6784.636399999999
$ws.Range("H129").Value = 726.5454999999999
$ws.Range("I129").Value = 563.8889
$ws.Range("J129").Value = 1458.5
$ws.Range("K129").Value = 1691.6667
$ws.Range("L129").Value = 4375.5
$ws.Range("M129").Value = 3308.3333
$ws.Range("N129").Value = -14375.5
$ws.Range("H137").Value = 2703.0166
$ws.Range("I137").Value = 1420.8
$ws.Range("J137").Value = 3985.2334
$ws.Range("K137").Value = 4262.4
$ws.Range("L137").Value = 11955.7002
$ws.Range("M137").Value = -1712.4
$ws.Range("N137").Value = -17055.7002
$ws.Range("H138").Value = 1025419.7
$ws.Range("I138").Value = 1646.1111
$ws.Range("J138").Value = 1284967.9
$ws.Range("K138").Value = 4938.3333
$ws.Range("L138").Value = 3854903.7
$ws.Range("M138").Value = 201.6666999999998
$ws.Range("N138").Value = -3865183.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10292.099
$ws.Range("I32").Value = 7676.0127
$ws.Range("J32").Value = 29318.182
$ws.Range("K32").Value = 7676.0127
$ws.Range("L32").Value = 29318.182
$ws.Range("M32").Value = -7389.0127
$ws.Range("N32").Value = -29892.182
$ws.Range("H61").Value = 10541.667
$ws.Range("I61").Value = 5259.35
$ws.Range("J61").Value = 36953.25
$ws.Range("K61").Value = 5259.35
$ws.Range("L61").Value = 36953.25
$ws.Range("M61").Value = -5047.35
$ws.Range("N61").Value = -37377.25
$ws.Range("H97").Value = 2057.6667
$ws.Range("I97").Value = 1323.8
$ws.Range("J97").Value = 2975
$ws.Range("K97").Value = 1323.8
$ws.Range("L97").Value = 2975
$ws.Range("M97").Value = -827.8
$ws.Range("N97").Value = -3967
$ws.Range("H101").Value = 49701.25
$ws.Range("J101").Value = 49701.25
$ws.Range("L101").Value = 49701.25
$ws.Range("N101").Value = -56191.25
$ws.Range("H122").Value = 8335466.5
$ws.Range("I122").Value = 2250
$ws.Range("K122").Value = 6750
$ws.Range("M122").Value = -4300
$ws.Range("H132").Value = 5325.512
$ws.Range("I132").Value = 1906.3
$ws.Range("K132").Value = 5718.9
$ws.Range("M132").Value = -3188.9
$ws.Range("H136").Value = 10541.667
$ws.Range("I136").Value = 5259.35
$ws.Range("J136").Value = 36953.25
$ws.Range("K136").Value = 15778.05
$ws.Range("L136").Value = 110859.75
$ws.Range("M136").Value = -13228.05
$ws.Range("N136").Value = -115959.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 255.29167
$ws.Range("I80").Value = 138.8
$ws.Range("K80").Value = 138.8
$ws.Range("M80").Value = 859.2
$ws.Range("H83").Value = 255.29167
$ws.Range("I83").Value = 138.8
$ws.Range("K83").Value = 694
$ws.Range("M83").Value = 4298
$ws.Range("H99").Value = 2320.889
$ws.Range("I99").Value = 2230
$ws.Range("J99").Value = 2378.7273
$ws.Range("K99").Value = 2230
$ws.Range("L99").Value = 2378.7273
$ws.Range("M99").Value = -732
$ws.Range("N99").Value = -5374.7273
$ws.Range("H132").Value = 62981.668
$ws.Range("J132").Value = 62981.668
$ws.Range("L132").Value = 62981.668
$ws.Range("N132").Value = -73101.66800000001
$ws.Range("H134").Value = 23762.912
$ws.Range("I134").Value = 2068.4443
$ws.Range("J134").Value = 1000014
$ws.Range("K134").Value = 6205.3329
$ws.Range("L134").Value = 3000042
$ws.Range("M134").Value = -3670.3329
$ws.Range("N134").Value = -3005112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5930.05
$ws.Range("I31").Value = 1044.1428
$ws.Range("J31").Value = 17330.5
$ws.Range("K31").Value = 1044.1428
$ws.Range("L31").Value = 17330.5
$ws.Range("M31").Value = -749.1428000000001
$ws.Range("N31").Value = -17920.5
$ws.Range("H34").Value = 5930.05
$ws.Range("I34").Value = 1044.1428
$ws.Range("J34").Value = 17330.5
$ws.Range("K34").Value = 1044.1428
$ws.Range("L34").Value = 17330.5
$ws.Range("M34").Value = -842.1428000000001
$ws.Range("N34").Value = -17734.5
$ws.Range("H58").Value = 1979195.5
$ws.Range("I58").Value = 3369220.5
$ws.Range("J58").Value = 3896.842
$ws.Range("K58").Value = 3369220.5
$ws.Range("L58").Value = 3896.842
$ws.Range("M58").Value = -3369017.5
$ws.Range("N58").Value = -4302.842000000001
$ws.Range("H132").Value = 2812.6
$ws.Range("I132").Value = 2660
$ws.Range("J132").Value = 3168.6667
$ws.Range("K132").Value = 7980
$ws.Range("L132").Value = 9506.000100000001
$ws.Range("M132").Value = -5450
$ws.Range("N132").Value = -14566.0001
$ws.Range("H136").Value = 1979195.5
$ws.Range("I136").Value = 3369220.5
$ws.Range("J136").Value = 3896.842
$ws.Range("K136").Value = 10107661.5
$ws.Range("L136").Value = 11690.526
$ws.Range("M136").Value = -10105111.5
$ws.Range("N136").Value = -16790.526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3475171.2
$ws.Range("I5").Value = 562.2632
$ws.Range("J5").Value = 16678686
$ws.Range("K5").Value = 1686.7896
$ws.Range("L5").Value = 50036058
$ws.Range("M5").Value = -1574.7896
$ws.Range("N5").Value = -50036282
$ws.Range("H129").Value = 2644.1538
$ws.Range("J129").Value = 1733
$ws.Range("L129").Value = 5199
$ws.Range("N129").Value = -15199
$ws.Range("H135").Value = 3475171.2
$ws.Range("I135").Value = 562.2632
$ws.Range("J135").Value = 16678686
$ws.Range("K135").Value = 5060.3688
$ws.Range("L135").Value = 150108174
$ws.Range("M135").Value = -2525.3688
$ws.Range("N135").Value = -150113244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12251.625
$ws.Range("I122").Value = 10752.333
$ws.Range("J122").Value = 16749.5
$ws.Range("K122").Value = 32256.999
$ws.Range("L122").Value = 50248.5
$ws.Range("M122").Value = -29806.999
$ws.Range("N122").Value = -55148.5
$ws.Range("H126").Value = 2825.6843
$ws.Range("I126").Value = 1724
$ws.Range("J126").Value = 4714.2856
$ws.Range("K126").Value = 5172
$ws.Range("L126").Value = 14142.8568
$ws.Range("M126").Value = -2702
$ws.Range("N126").Value = -19082.8568
$ws.Range("H132").Value = 43497
$ws.Range("I132").Value = 64526.688
$ws.Range("J132").Value = 12908.363
$ws.Range("K132").Value = 193580.064
$ws.Range("L132").Value = 38725.089
$ws.Range("M132").Value = -191050.064
$ws.Range("N132").Value = -43785.089
$ws.Range("H133").Value = 34666.668
$ws.Range("J133").Value = 34666.668
$ws.Range("L133").Value = 34666.668
$ws.Range("N133").Value = -44786.668
$ws.Range("H139").Value = 42326
$ws.Range("J139").Value = 42326
$ws.Range("L139").Value = 42326
$ws.Range("N139").Value = -52606

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 5500
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10340
$ws.Range("H40").Value = 4318
$ws.Range("I40").Value = 3995.8262
$ws.Range("K40").Value = 3995.8262
$ws.Range("M40").Value = -3859.8262
$ws.Range("H94").Value = 48995
$ws.Range("J94").Value = 48995
$ws.Range("L94").Value = 48995
$ws.Range("N94").Value = -50347
$ws.Range("H132").Value = 5753.65
$ws.Range("I132").Value = 6701.846
$ws.Range("J132").Value = 3992.7144
$ws.Range("K132").Value = 20105.538
$ws.Range("L132").Value = 11978.1432
$ws.Range("M132").Value = -17575.538
$ws.Range("N132").Value = -17038.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1685.7142
$ws.Range("I126").Value = 1685.7142
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5057.142599999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2587.142599999999
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2437.5557
$ws.Range("I132").Value = 1467.1333
$ws.Range("J132").Value = 7289.6665
$ws.Range("K132").Value = 4401.3999
$ws.Range("L132").Value = 21868.9995
$ws.Range("M132").Value = -1871.3999
$ws.Range("N132").Value = -26928.9995
$ws.Range("H138").Value = 41488.625
$ws.Range("J138").Value = 41488.625
$ws.Range("L138").Value = 41488.625
$ws.Range("N138").Value = -51768.625

Write-Host "Applied all Pandaemonium Profits updates."